$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.391.86'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '2.105.56'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  -0.96%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '343.48'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("E6").Value = '  -0.93%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5314'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.15%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.4441'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.52%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '54.86'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.66%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.09415'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.98%  '
$ws.Range("E11").Value = '  +0.89%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '24.87'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.51%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '8.580'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.41%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.929'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.74%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.071.94'
$ws.Range("E15").Value = '  -1.60%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '101.85'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.42%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001162'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("E19").Value = '  +2.65%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.06692'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.55%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.342'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.70%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.89%  '
$ws.Range("D23").Value = '30.427.27'
$ws.Range("E23").Value = '  +2.27%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '12.57'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.97%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.311'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.34%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '21.93'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.10%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '162.76'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.09%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.531'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.74%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.801'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +7.95%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '133.86'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.81%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.150'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.674'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +3.28%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.1056'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("E34").Value = '  +1.49%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.849'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("E36").Value = '  +0.01%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.02649'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.06%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.06809'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.60%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '12.72'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.64%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.7042'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.70%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.347'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.27%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.2226'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("E43").Value = '  -1.45%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '14.52'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.67%  '
$ws.Range("E45").Value = '  +1.77%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.86%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.384'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +19.09%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '3.642'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.57%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.00000000351'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.27%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.222'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +8.66%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.220'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.19%  '
